$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for data rows 2-13
# from serial date 45183 (2023-09-14) to 45184 (2023-09-15).
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
